$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the data block (rows 3-21, columns A-E) one column to the right
#     (A->B, B->C, C->D, D->E, E->F), to make room for a new leading
#     "BalanceSheetSide" column. Rows 1-2 (the Template / Item_type header
#     rows) are left untouched. Work right-to-left so we never clobber a
#     source cell before it has been read.
for ($r = 3; $r -le 21; $r++) {
    for ($c = 5; $c -ge 1; $c--) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r, $c + 1)
        $dst.Value = $src.Value()
        $dst.NumberFormat = $src.NumberFormat()
    }
    # Column A no longer holds the old content - clear it (it will be
    # repopulated below for the rows that need the new BalanceSheetSide
    # label) and reset it back to a plain/general cell.
    $a = $ws.Cells.Item($r, 1)
    $a.ClearContents()
    $a.NumberFormat = "General"
}

# --- New leading column content
$ws.Cells.Item(7, 1).Value = "BalanceSheetSide"
$ws.Cells.Item(8, 1).Value = "Assets"
$ws.Cells.Item(9, 1).Value = "Assets"
$ws.Cells.Item(10, 1).Value = "Assets"

# --- Column widths: column A is new (sized like the old "Item_type" column),
#     column B keeps the old date-column width, column C keeps the old
#     percent-header column width.
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Columns.Item(2).ColumnWidth = 9.72
$ws.Columns.Item(3).ColumnWidth = 19.39

# --- Selection follows the shifted "counter item type" row value cell
$ws.Range("F6").Select()
